{"js": "// The \"COMPETENCES TECHNIQUES\" skills block is made up of 8 one-line\n// paragraphs. The resume_data source got converted from a Python dict to\n// JSON, which changed the iteration/key order used to emit these lines.\n// The paragraphs themselves (and their formatting) stay put - only the\n// text content slotted into each paragraph changes, in this exact order:\n//   Web, Langages, Bases de donn\u00e9es, Soft_Skills, Autres, Visualisation, ML/AI, MLOps\n// becomes\n//   Langages, Soft_Skills, Visualisation, MLOps, Web, Autres, ML/AI, Bases de donn\u00e9es\n\nconst oldOrder = [\n  \"Web : client\",\n  \"Langages : r, python, matlab, c, c++\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Soft_Skills : dashboards\",\n  \"Autres : marketing,  prodigi \u2013 agile ,  activeviam ,  data ,  r , vulgarise, webinar\",\n  \"Visualisation : tableau\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : ux, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n];\n\nconst newOrder = [\n  \"Langages : r, python, matlab, c, c++\",\n  \"Soft_Skills : dashboards\",\n  \"Visualisation : tableau\",\n  \"MLOps : ux, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Web : client\",\n  \"Autres : marketing,  prodigi \u2013 agile ,  activeviam ,  data ,  r , vulgarise, webinar\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the contiguous run of paragraphs that holds the skills block by\n// matching the first (and longest, to be safe) expected sequence.\nconst items = paragraphs.items;\nlet startIndex = -1;\nfor (let i = 0; i + oldOrder.length <= items.length; i++) {\n  let matches = true;\n  for (let j = 0; j < oldOrder.length; j++) {\n    if (items[i + j].text !== oldOrder[j]) {\n      matches = false;\n      break;\n    }\n  }\n  if (matches) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not locate the COMPETENCES TECHNIQUES skills paragraphs.\");\n}\n\nfor (let j = 0; j < newOrder.length; j++) {\n  items[startIndex + j].insertText(newOrder[j], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The \"COMPETENCES TECHNIQUES\" skills block is made up of 8 one-line\n# paragraphs. The resume_data source got converted from a Python dict to\n# JSON, which changed the iteration/key order used to emit these lines.\n# The paragraphs themselves (and their formatting) stay put - only the\n# text content slotted into each paragraph changes, in this exact order:\n#   Web, Langages, Bases de donn\u00e9es, Soft_Skills, Autres, Visualisation, ML/AI, MLOps\n# becomes\n#   Langages, Soft_Skills, Visualisation, MLOps, Web, Autres, ML/AI, Bases de donn\u00e9es\n\n$d = $word.ActiveDocument\n\n$oldOrder = @(\n  \"Web : client\",\n  \"Langages : r, python, matlab, c, c++\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Soft_Skills : dashboards\",\n  \"Autres : marketing,  prodigi \u2013 agile ,  activeviam ,  data ,  r , vulgarise, webinar\",\n  \"Visualisation : tableau\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : ux, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n)\n\n$newOrder = @(\n  \"Langages : r, python, matlab, c, c++\",\n  \"Soft_Skills : dashboards\",\n  \"Visualisation : tableau\",\n  \"MLOps : ux, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Web : client\",\n  \"Autres : marketing,  prodigi \u2013 agile ,  activeviam ,  data ,  r , vulgarise, webinar\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n)\n\n$n = $oldOrder.Length\n$count = $d.Paragraphs.Count\n\n# Locate the contiguous run of paragraphs holding the skills block by\n# matching the expected text sequence (trim the trailing paragraph mark\n# that Range.Text includes).\n$startIndex = -1\nfor ($i = 1; $i -le ($count - $n + 1); $i++) {\n    $allMatch = $true\n    for ($j = 0; $j -lt $n; $j++) {\n        $t = $d.Paragraphs.Item($i + $j).Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -ne $oldOrder[$j]) {\n            $allMatch = $false\n            break\n        }\n    }\n    if ($allMatch) {\n        $startIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1) {\n    throw \"Could not locate the COMPETENCES TECHNIQUES skills paragraphs.\"\n}\n\nfor ($j = 0; $j -lt $n; $j++) {\n    $d.Paragraphs.Item($startIndex + $j).Range.Text = $newOrder[$j]\n}\n"}
